# database schema.xlsx -- add a new "Sheet1" wireframe/component sheet
# continuing the App.vue component-schema diagram (Homepage, login/register,
# product search, individual product, view cart pages), and make it the
# active sheet.

$wb = $excel.ActiveWorkbook

# --- add the new worksheet at the end of the tab strip --------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $last)

# --- column widths (matches the componentSchema-style wireframe layout) ---
$ws.Columns("A:D").ColumnWidth = 9.140625
$ws.Columns("E:E").ColumnWidth = 4
$ws.Columns("F:F").ColumnWidth = 31.140625
$ws.Columns("G:G").ColumnWidth = 4
$ws.Columns("H:H").ColumnWidth = 15.5703125
$ws.Columns("I:I").ColumnWidth = 4
$ws.Columns("J:J").ColumnWidth = 9.140625
$ws.Columns("K:K").ColumnWidth = 4

# --- cell content -----------------------------------------------------
$ws.Range("D3").Value = "App.vue"
$ws.Range("F10").Value = "Navbar"

$ws.Range("F5").Value = "Homepage"
$ws.Range("H5").Value = "banner"
$ws.Range("H7").Value = "card small"
$ws.Range("F26").Value = "product search Page"
$ws.Range("F14").Value = "user login page"
$ws.Range("F18").Value = "user register page"
$ws.Range("F35").Value = "view cart page"
$ws.Range("F32").Value = "individual product page"
$ws.Range("F39").Value = "item you sell"
$ws.Range("F28").Value = "per seller item sold"
$ws.Range("H26").Value = "cardSmall"
$ws.Range("H35").Value = "transactionForm"
$ws.Range("H32").Value = "product info"
$ws.Range("H14").Value = "user form"

$ws.Range("H18").Value = "user form"
$ws.Range("H28").Value = "cardSmall"
$ws.Range("H39").Value = "cardSmall"

# --- center all populated cells (matches the sheet's base style) ----------
$labelCells = @("D3","F5","H5","H7","F10","F14","H14","F18","H18","F26","H26","F28","H28","F32","H32","F35","H35","F39","H39")
foreach ($addr in $labelCells) {
    $ws.Range($addr).HorizontalAlignment = -4108
    $ws.Range($addr).VerticalAlignment = -4108
}

# --- black divider bars between the label/value column pairs --------------
$ws.Range("E39").Interior.ThemeColor = 1
$ws.Range("G39").Interior.ThemeColor = 1
$ws.Range("I39").Interior.ThemeColor = 1
$ws.Range("K39").Interior.ThemeColor = 1

# --- this new sheet becomes the active / selected tab ----------------------
$selected = $ws.Range("H2").Select()
$activated = $ws.Activate()
